$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Objective paragraph: drop the trailing "...that offers competitive pay
#    and benefits for me and my family" clause so the sentence now reads
#    "...with a stable and growing company." instead of continuing on.
# ---------------------------------------------------------------------------
$ok1 = $d.Content.Find.Execute(
    " that offers competitive pay and benefits for me and my family",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 2)
if (-not $ok1) {
    throw "Objective text not found"
}

# ---------------------------------------------------------------------------
# 2) Coursework paragraph: reorder / rewrite the list of classes.
# ---------------------------------------------------------------------------
$oldCoursework = "Programming 1, Programming 2, Programming in C, Circuits 1, Digital Logic, Microprocessor Applications, Electrical Junior Design, Electrical Senior Design, Software Senior Design, Database Systems 1 & 2, Data Structures, Digital Design, Intro. Software Engineering, Operating Systems"
$newCoursework = "Programming 1 & 2, Intro. Software Engineering, Database Systems 1 & 2, Digital Design, Circuits 1, Digital Logic, Microprocessor Applications, Electrical Junior Design, Electrical Senior Design, Software Senior Design, Data Structures,  Operating Systems"
$ok2 = $d.Content.Find.Execute(
    $oldCoursework, $true, $false, $false, $false, $false, $true, 1, $false,
    $newCoursework, 2)
if (-not $ok2) {
    throw "Coursework text not found"
}

# ---------------------------------------------------------------------------
# 3) Move the "_GoBack" bookmark from the Projects bullet ("for Graduate
#    students...") to the blank paragraph that follows the Coursework line.
# ---------------------------------------------------------------------------
try {
    $goBack = $d.Bookmarks("_GoBack")
    $goBack.Delete()
} catch {
    # No existing _GoBack bookmark -- nothing to remove.
}

# Locate the Coursework paragraph by its (now-updated) content, then grab the
# empty paragraph that immediately follows it -- that's the new bookmark home.
$courseRng = $d.Content
$ok3 = $courseRng.Find.Execute("Operating Systems")
if (-not $ok3) {
    throw "Coursework paragraph not found for bookmark placement"
}
$coursePara = $courseRng.Paragraphs(1)
$blankPara = $coursePara.Next()
$bmRange = $blankPara.Range
$bmRange.Collapse(1)
$d.Bookmarks.Add("_GoBack", $bmRange)

Write-Output "Objective shortened: $ok1"
Write-Output "Coursework rewritten: $ok2"
Write-Output "_GoBack bookmark relocated"
